# Rrxcell::Sheet#address / Rrxcell::Book#sheet_names/#sheets/#sheet support work
# added a small "address lookup" demo block to Sheet1: three 2-column groups
# (AA, AB/AC) used to exercise the new A1-format addressing code paths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New narrow "spacer" columns (F through Z) - matches the worksheet's
# <col min="6" max="26" width="3.33203125" customWidth="1"/> layout used to
# keep the new AA/AB/AC demo columns visually separated from the existing
# A:E data block.
$ws.Range("F1:Z1").EntireColumn.ColumnWidth = 2.57

# New sample data used by the address-lookup examples - written column by
# column (AA first, then AB/AC row by row) which is the order the new
# shared-string entries show up in the saved workbook.
$ws.Range("AA1").Value = "Sheet1!AA1"
$ws.Range("AA2").Value = "Sheet1!AA2"
$ws.Range("AA3").Value = "Sheet1!AA3"

$ws.Range("AB1").Value = "Sheet1!AB1"
$ws.Range("AC1").Value = "Sheet1!AC1"

$ws.Range("AB2").Value = "Sheet1!AB2"
$ws.Range("AC2").Value = "Sheet1!AC2"

$ws.Range("AB3").Value = "Sheet1!AB3"
$ws.Range("AC3").Value = "Sheet1!AC3"

# Move the active selection over to the new block.
[void]$ws.Range("Q7").Select()

# Page setup tweak that came along with this commit's save.
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
